$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-18.82206973900809, 1.969469173410194, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 1.163695777329076, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, 1.959434857934247, 1.670429928330072, -18.82206973900809, 3.413466861024522, -18.82206973900809, -18.82206973900809, -18.82206973900809, 0.726761574505961, -18.82206973900809),
    @(-18.82206973900809, 1.682156679027108, -18.82206973900809, -18.82206973900809, -18.82206973900809, 2.878085888742318, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(2.452145163541375, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, 1.820853645577504, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(3.860731367538047, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 1.762082610653426, -18.82206973900809, 2.174989113438325),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, 2.918436139011992, -18.82206973900809, 2.743776576360386, -18.82206973900809, -18.82206973900809, -18.82206973900809, 2.088910351777756),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, 2.488266231823164, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 1.587612267623396, 1.822966786913098),
    @(-18.82206973900809, -18.82206973900809, 1.525864079884549, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 1.840526495781883),
    @(-18.82206973900809, -18.82206973900809, 1.718971087322718, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 2.00625969795171, -18.82206973900809),
    @(-18.82206973900809, 2.131320357783188, 1.869252076137157, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 2.052858343225132, 2.548835261920447, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809, 2.020370383300361, 2.451678631163457, -18.82206973900809),
    @(-18.82206973900809, -18.82206973900809, 2.055421854494116, -18.82206973900809, -18.82206973900809, -18.82206973900809, 4.321925137617091, 1.897784455980642, -18.82206973900809, -18.82206973900809),
    @(-18.82206973900809, 1.044069679625614, 1.50344324493831, -18.82206973900809, 3.224179349157219, -18.82206973900809, -18.82206973900809, 1.284059929166576, -18.82206973900809, 2.039275438878529),
    @(-18.82206973900809, 1.346157980433512, -18.82206973900809, 1.720615484892776, -18.82206973900809, 2.572891630339132, -18.82206973900809, -18.82206973900809, -18.82206973900809, -18.82206973900809)
)

$data = New-Object 'object[,]' 20,10
for ($r = 0; $r -lt $values.Count; $r++) {
    for ($c = 0; $c -lt $values[$r].Count; $c++) {
        $data[$r, $c] = $values[$r][$c]
    }
}

$ws.Range("B2:K21").Value = $data
